# Update the single sample data row (ORDEN / NUM_IMEI / NUM_SIMCARD) on the
# "Orden" sheet with the new test fixture values, refresh the saved
# selection, and drop the now-unused trailing blank (but formatted) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophes force these numeric-looking strings to stay text
# (matching the workbook's existing quotePrefix text styles) instead of
# being coerced to numbers (which would also mangle the 19-digit SIM
# number into scientific notation).
$ws.Range("B2").Value = "'1005069"
$ws.Range("D2").Value = "'8954080008100062856"
$ws.Range("C2").Value = "'123600000007459"

$ws.Range("C8").Select()

$ws.Rows("55:71").Delete()
